$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 3: ASP label stays the same text, no change needed ---

# --- Row 4: ITN proportion assumption 0.71 -> 0.7 ---
$ws.Range("B4").Value = 0.7

# --- Row 6: clear out the old "ITN proportion" row (label + value) ---
$ws.Range("A6").ClearContents()
$ws.Range("C6").ClearContents()

# --- Row 8: Cost of distribution % 0.29 -> 0.3, and formula now refers to B4/B8 ---
$ws.Range("B8").Value = 0.3
$ws.Range("C8").Formula = "=(C4/B4)*B8"

# --- A16 / A17 label text updates ---
# A16 previously had rich (bold) formatting - the new text is applied as a
# plain replacement (matches how the source file now stores it).
$ws.Range("A16").Value = "Total physical degradation loss" + [char]10 + "(Insecticide efficacy + wear and tear)"

# A17 keeps its original rich-text run structure (bold "Total lost value" +
# normal remainder); only trailing " loss" is dropped from the final run, so
# edit it in place through Characters() to preserve the two runs.
$a17 = $ws.Range("A17")
$a17.Characters(18, 91).Text = "(LLIN lost + not used every night + minimal insecticide efficacy loss + wear and tear)"
$a17.Characters(1, 16).Font.Bold = $true
$a17.Characters(17, 92).Font.Bold = $false

# --- A41 label text update (also now needs wrap text since multi-line) ---
$ws.Range("A41").Value = "Improve physical degradation by " + [char]10 + "(insecticide efficacy + wear and tear)"
$ws.Range("A41").WrapText = $true

# --- View state: selection moved to C6, scrolled back to top-left ---
$ws.Range("C6").Select()
